$wb = $excel.ActiveWorkbook

# "Latest Handback DateTime" for row 2 (file 5f266658-...) was refreshed when
# the handback report was regenerated.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-10-25 03:08:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-25 03:09:14"
